# NSMB - 5-2 - 6 frames saved in 5-2 thanks to terrotim.
# Updates the splits table on sheet "V4" (sheet1): fills in newly-recorded
# frame counts (columns B/C), refreshes the dependent diff formulas
# (columns D/F), clears the now-stale "G"/"H" columns for rows that no
# longer have a recorded comparison point, relabels several checkpoint
# rows with the new, more precise checkpoint names, and removes the two
# trailing rows that are no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the two trailing rows (old rows 56 & 57) -----------------------
# Everything below row 55 collapses up by two rows.
$ws.Rows("56:57").Delete()

# --- Column B (new frame counts recorded for this run) ---------------------
$ws.Range("B39").Value = 15290
$ws.Range("B40").Value = 15631
$ws.Range("B44").Value = 16363
$ws.Range("B45").Value = 16557
$ws.Range("B46").Value = 16685
$ws.Range("B47").Value = 16859
$ws.Range("B48").Value = 17061
$ws.Range("B49").Value = 17390
$ws.Range("B50").Value = 17549
$ws.Range("B51").Value = 17569
$ws.Range("B52").Value = 17640

# Rows 54 & 55 no longer carry a recorded "B" frame count.
$ws.Range("B54").Value = ""
$ws.Range("B55").Value = ""

# --- Column C (new frame counts for the comparison column) -----------------
$ws.Range("C36").Value = 15235
$ws.Range("C37").Value = 15489
$ws.Range("C38").Value = 17048
$ws.Range("C39").Value = 17562
$ws.Range("C40").Value = 18105
$ws.Range("C42").Value = 18403
$ws.Range("C43").Value = 18789
$ws.Range("C44").Value = 19022
$ws.Range("C45").Value = 19216
$ws.Range("C46").Value = 19344
$ws.Range("C47").Value = 19523
$ws.Range("C48").Value = 19725
$ws.Range("C49").Value = 20054
$ws.Range("C50").Value = 20197
$ws.Range("C51").Value = 20233
$ws.Range("C52").Value = 20304

# --- Column E (rows 53-55 now carry the values that used to live on the
#     deleted rows 56/57) ----------------------------------------------------
$ws.Range("E53").Value = 20600
$ws.Range("E54").Value = 21114
$ws.Range("E55").Value = 21615

# --- Re-stripe the shared "diff" formulas so they cover the new row ranges -
$ws.Range("D36:D55").Formula = "=IF(B36 >  0,C36-B36, 0)"
$ws.Range("F53:F55").Formula = "=IF(B53 >  0,E53-B53, 0)"

# --- Column A relabeling (checkpoints renamed / renumbered) ----------------
$ws.Range("A40").Value = "Map First move"
$ws.Range("A44").Value = "Checkpoint 9"
$ws.Range("A45").Value = "Checkpoitn Rail 460xxxxx"
$ws.Range("A46").Value = "Checkpoint 995"
$ws.Range("A47").Value = "Checkpoint 1500"
$ws.Range("A48").Value = "Checkpoint 2107/2106"
$ws.Range("A49").Value = "Checkpoint 3094/3093"
$ws.Range("A50").Value = "Checkpoint 3442"
$ws.Range("A51").Value = "Enter pipe"
$ws.Range("A52").Value = "Black screen"
$ws.Range("A53").Value = "Get flag"
$ws.Range("A54").Value = "Black screen"
$ws.Range("A55").Value = "1st Move"

# --- Columns G/H: clear the now-unused "end of level" comparison values ----
$clearRows = 42,43,45,46,47,48,49,50,51,52,53,54,55
foreach ($r in $clearRows) {
    $ws.Range("G$r").Value = ""
    $ws.Range("H$r").Value = ""
}

# --- Scroll/selection bookkeeping (matches the trimmed row count) ----------
$ws.Range("B52").Select()
$ws.Application.ActiveWindow.ScrollRow = 34
